$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column G with the same visual formatting as column F, then give it
# the requested header text ("Translation Error") in row 13, matching the
# existing header row style used by B13:F13.
$ws.Range("F13").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("G13").Value = "Translation Error"

# Size the new column similarly to the reference workbook (~30.81 characters).
$ws.Columns("G").ColumnWidth = 29.917

# Update the visible scroll position / active selection like the source edit.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B12").Select()
